$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1575:C1575").Value = "Wardah Everyday Cheek Liptint"
$ws.Range("B1576:C1576").Value = "Wardah Exclusive Matte Lip Cream"
$ws.Range("B1577:C1577").Value = "Wardah Eyexpert Eye Lip Make Up Remover"
$ws.Range("B1578:C1578").Value = "Wardah Intense Matte Lipstick"
$ws.Range("B1579:C1579").Value = "Wardah Lipstick Longlasting"
$ws.Range("B1580:C1580").Value = "Wardah Longlasting Lipstick"
$ws.Range("B1581:C1581").Value = "Wardah Eyexpert Optimum Hi Black Liner"
$ws.Range("B1582:C1582").Value = "Wardah Eyexpert Aqua Lash Mascara"
$ws.Range("B1583:C1583").Value = "Wardah Eyexpert Eye Shadow Classic"
$ws.Range("B1584:C1584").Value = "Wardah Eyexpert The Volume Expert Mascara"
$ws.Range("B1585:C1585").Value = "Wardah Renew You Anti Aging Eye Cream"
$ws.Range("B1586:C1586").Value = "Wardah Eyexpert Perfectcurl Mascara"
$ws.Range("B1587:C1587").Value = "Wardah Eyexpert Staylast Liquid Eyeliner"
$ws.Range("B1588:C1588").Value = "Wardah Eye Shadow"
$ws.Range("B1589:C1589").Value = "Wardah Eyexpert Eye Shadow"
$ws.Range("B1590:C1590").Value = "Wardah Acnederm Face Powder"
$ws.Range("B1591:C1591").Value = "Wardah Blush On"
$ws.Range("B1592:C1592").Value = "Wardah Exclusive Liquid Foundation"
$ws.Range("B1593:C1593").Value = "Wardah Exclusive Two Way Cake"
$ws.Range("B1594:C1594").Value = "Wardah Everyday Luminous Compact Powder"
$ws.Range("B1595:C1595").Value = "Wardah Everyday Luminous Face Powder"
$ws.Range("B1596:C1596").Value = "Wardah Refill Lightening Powder Foundation Light Feel"
$ws.Range("B1597:C1597").Value = "Wardah Lightening Powder Foundation Light Feel"
$ws.Range("B1598:C1598").Value = "Wardah Refill Exclusive Two Way Cake"
$ws.Range("B1599:C1599").Value = "Wardah Lightening Powder Foundation"
$ws.Range("B1600:C1600").Value = "Wardah Refill Lightening Powder Foundation"
$ws.Range("B1601:C1601").Value = "Wardah Everyday Luminous Liquid Foundation"
$ws.Range("B1602:C1602").Value = "Wardah Refill Everyday Luminous Two Way Cake"
$ws.Range("B1603:C1603").Value = "Wardah Everyday Shine Free Bb Loose Powder"
$ws.Range("B1604:C1604").Value = "Wardah Lightening Matte Powder"
$ws.Range("B1605:C1605").Value = "Wardah Lightening Bb Cake Powder"
$ws.Range("B1606:C1606").Value = "Wardah Everyday Luminous Two Way Cake"
$ws.Range("B1607:C1607").Value = "Wardah Instaperfect Porefection Skin Primer"
$ws.Range("B1608:C1608").Value = "Wardah Instaperfect Mineralight Matte Bb Cushion"
$ws.Range("B1609:C1609").Value = "Wardah Instaperfect Spotlight Chromatic Eye Palette"
$ws.Range("B1610:C1610").Value = "Wardah Instaperfect Geniustwist Matic Contour Brow Brushed"
$ws.Range("B1611:C1611").Value = "Wardah Instaperfect Hypergetic Precise Black Liner"
$ws.Range("B1612:C1612").Value = "Wardah Instaperfect Mattesetter Lip Matte Paint"
$ws.Range("B1613:C1613").Value = "Wardah Instaperfect Quick Fix Cover Correct Concealer"
$ws.Range("B1614:C1614").Value = "Wardah Instaperfect Dynamatic Microsmooth Liner"
$ws.Range("B1615:C1615").Value = "Wardah Refill Instaperfect Mineralight Matte Bb Cushion"
$ws.Range("B1616:C1616").Value = "Wardah Instaperfect Gloss Chic Lip Crayon"
$ws.Range("B1617:C1617").Value = "Wardah Instaperfect Mattetitude Matte Stain Lipstick"
$ws.Range("B1618:C1618").Value = "Wardah Instaperfect Mattecentric Lip Crayon"
$ws.Range("B1619:C1619").Value = "Wardah Instaperfect Matte Fit Powder Foundation"
$ws.Range("B1620:C1620").Value = "Wardah Refill Instaperfect Matte Fit Powder Foundation"
$ws.Range("B1621:C1621").Value = "Wardah Colorfit Velvet Matte Lip Mousse"
$ws.Range("B1622:C1622").Value = "Emina Magic Potion"
$ws.Range("B1623:C1623").Value = "Emina Glossy Stain"
$ws.Range("B1624:C1624").Value = "Emina Cheek Lit Pressed Blush"
$ws.Range("B1625:C1625").Value = "Emina Cheek Lit Cream Blush"
$ws.Range("B1626:C1626").Value = "Emina Beauty Bliss Bb Cream"
$ws.Range("B1627:C1627").Value = "Emina Bare With Me Mineral Cushion"
$ws.Range("B1628:C1628").Value = "Emina Bright Stuff Loose Powder"
$ws.Range("B1629:C1629").Value = "Emina Pore Ranger"
$ws.Range("B1630:C1630").Value = "Emina Daily Matte Loose Powder"
$ws.Range("B1631:C1631").Value = "Emina Pop Rouge Pressed Eye Shadow"
$ws.Range("B1632:C1632").Value = "Emina City Chic Cc Cake"
$ws.Range("B1633:C1633").Value = "Emina Star Lash Aqua Mascara"
$ws.Range("B1634:C1634").Value = "Emina Agent Of Brow"
$ws.Range("B1635:C1635").Value = "Emina Eye Do Crayon Pour Les Yeux"
$ws.Range("B1636:C1636").Value = "Emina Top Secret Eyebrow"
$ws.Range("B1637:C1637").Value = "Emina Double Agent Eyebrow"
$ws.Range("B1638:C1638").Value = "Emina Creamatte Metallic Edition"
$ws.Range("B1639:C1639").Value = "Emina Soulmatte Lipstick"
$ws.Range("B1640:C1640").Value = "Emina Sugar Rush Lipstick"
$ws.Range("B1641:C1641").Value = "Emina Creamytint"
$ws.Range("B1642:C1642").Value = "Emina Creamatte"
$ws.Range("B1643:C1643").Value = "Emina Lip Cushion"
$ws.Range("B1644:C1644").Value = "Emina Tinted Lipbalm"
$ws.Range("B1645:C1645").Value = "Emina Liquid Lip Shine"
$ws.Range("B1646:C1646").Value = "Emina Smoochies Lipbalm"
$ws.Range("B1647:C1647").Value = "Sariayu Lipstick"
$ws.Range("B1648:C1648").Value = "Sariayu Color Trend 2019 Hydra Lip Tint"
$ws.Range("B1649:C1649").Value = "Sariayu Color Trend 2019 Lite Lip Cream"
$ws.Range("B1650:C1650").Value = "Sariayu Lip Care"
$ws.Range("B1651").Value = "Sariayu Color Trend 2016 Duo Lip Color"
$ws.Range("C1651").Value = "Sariayu Color Trend 2017 Duo Lip Color"
$ws.Range("B1652").Value = "Sariayu Color Trend 2017 Duo Lip Color"
$ws.Range("C1652").Value = "Sariayu Color Trend 2016 Duo Lip Color"
$ws.Range("B1653:C1653").Value = "Sariayu Color Trend 2018 Lip Cream"
$ws.Range("B1654:C1654").Value = "Sariayu Color Trend 2018 Lip Metallic"
$ws.Range("B1655:C1655").Value = "Sariayu Lip Colour Matte"
$ws.Range("B1656:C1656").Value = "Sariayu Ct 19 Lite Lip Cream"
$ws.Range("B1657:C1657").Value = "Sariayu Trend 2017 Duo Lip Color"
$ws.Range("B1658").Value = "Sariayu Lip Color Matte"
$ws.Range("C1658").Value = "Sariayu Lip Colour Matte"
$ws.Range("B1659:C1659").Value = "Sariayu Two Way Cake"
$ws.Range("B1660:C1660").Value = "Sariayu Moisturizer"
$ws.Range("B1661:C1661").Value = "Sariayu Color Trend 2020 Lip Cheek"
$ws.Range("B1662:C1662").Value = "Sariayu Two Way Cake Energizing Aromatic Refill"
$ws.Range("B1663:C1663").Value = "Sariayu Alas Bedak Energizing Aromatic"
$ws.Range("B1664:C1664").Value = "Sariayu Color Trend 2020 Cheek Palette"
$ws.Range("B1665:C1665").Value = "Sariayu Loose Powder"
$ws.Range("B1666:C1666").Value = "Sariayu Creamy Foundation"
$ws.Range("B1667:C1667").Value = "Sariayu Compact Powder Spf 15"
$ws.Range("B1668:C1668").Value = "Sariayu Compact Powder"
$ws.Range("B1669:C1669").Value = "Sariayu Alas Bedak"
$ws.Range("B1670:C1670").Value = "Sariayu Bedak Jerawat Energizing Aromatic"
$ws.Range("B1671:C1671").Value = "Sariayu Refill Two Way Cake"
$ws.Range("B1672:C1672").Value = "Sariayu Blush On"
$ws.Range("B1673:C1673").Value = "Sariayu Color Trend 2015 Eyeliner Pencil"
$ws.Range("B1674:C1674").Value = "Sariayu Color Trend 2020 Eye Makeup Kit"
$ws.Range("B1675:C1675").Value = "Sariayu Pensil Alis Pro"
$ws.Range("B1676").Value = "Sariayu Color Trend 2019 Eye Shadow"
$ws.Range("C1676").Value = "Sariayu Color Trend 2016 Eye Shadow"
$ws.Range("B1677:C1677").Value = "Sariayu Color Trend 2015 Mascara"
$ws.Range("B1678").Value = "Sariayu Color Trend 16 Eye Shadow Kit"
$ws.Range("C1678").Value = "Sariayu Color Trend 18 Eye Shadow Kit"
$ws.Range("B1679").Value = "Sariayu Color Trend 18 Eye Shadow Kit"
$ws.Range("C1679").Value = "Sariayu Color Trend 16 Eye Shadow Kit"
$ws.Range("B1680:C1680").Value = "Sariayu Pensil Alis"
$ws.Range("B1681:C1681").Value = "Sariayu Trend Warna 2011 Moistpome Eye Shadow"
$ws.Range("B1682:C1682").Value = "Sariayu Color Trend Warna 2011 Moistpome Eye Shadow Palette"
$ws.Range("B1683:C1683").Value = "Sariayu Duo Eye Make Up"
$ws.Range("B1684:C1684").Value = "Sariayu Color Trend 2017 Liquid Eye Shadow"
$ws.Range("B1685").Value = "Sariayu Color Trend 2016 Eye Shadow"
$ws.Range("C1685").Value = "Sariayu Color Trend 2019 Eye Shadow"
$ws.Range("B1686:C1686").Value = "Sariayu Color Trend 2017 Eye Shadow Kit"
$ws.Range("B1687:C1687").Value = "Sariayu Color Trend 16 Eyeshadow"
